$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the stray empty trailing cells that were on the old last row (50)
$ws.Range("I50").ClearContents()
$ws.Range("K50:R50").ClearContents()

# Step 2: append new rows 51-72
# Row 51
$ws.Range("A51").Value = '4AF08793'
$ws.Range("B51").Value = 'SOL DE JANEIRO CHERIOSA 62 BODY MIST 240 ML'
$ws.Range("C51").Value = 'VARIOS'
$ws.Range("D51").Value = 'Tiene PT'
$ws.Range("E51").Value = 'Tiene ES'
$ws.Range("F51").Value = 'Tiene IT'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '1'
$ws.Range("H51").Value = 'UND'
$ws.Range("J51").Value = 'Solo Revisión'

# Row 52
$ws.Range("A52").Value = '4AF08794'
$ws.Range("B52").Value = 'SOL DE JANEIRO CHERIOSA SET BODY MIST'
$ws.Range("C52").Value = 'VARIOS'
$ws.Range("D52").Value = 'Tiene PT'
$ws.Range("E52").Value = 'Tiene ES'
$ws.Range("F52").Value = 'Tiene IT'
$ws.Range("G52").NumberFormat = "@"
$ws.Range("G52").Value = '1'
$ws.Range("H52").Value = 'UND'
$ws.Range("J52").Value = 'Solo Revisión'

# Row 53
$ws.Range("A53").Value = '4AF08795'
$ws.Range("B53").Value = 'SOL DE JANEIRO CHERIOSA 40 BODY MIST 90 ML'
$ws.Range("C53").Value = 'VARIOS'
$ws.Range("D53").Value = 'Tiene PT'
$ws.Range("E53").Value = 'Tiene ES'
$ws.Range("F53").Value = 'Tiene IT'
$ws.Range("G53").NumberFormat = "@"
$ws.Range("G53").Value = '1'
$ws.Range("H53").Value = 'UND'
$ws.Range("J53").Value = 'Solo Revisión'

# Row 54
$ws.Range("A54").Value = '4AF08796'
$ws.Range("B54").Value = 'SOL DE JANEIRO CHERIOSA 59 BODY MIST 90 ML'
$ws.Range("C54").Value = 'VARIOS'
$ws.Range("D54").Value = 'Tiene PT'
$ws.Range("E54").Value = 'Tiene ES'
$ws.Range("F54").Value = 'Tiene IT'
$ws.Range("G54").NumberFormat = "@"
$ws.Range("G54").Value = '1'
$ws.Range("H54").Value = 'UND'
$ws.Range("J54").Value = 'Solo Revisión'

# Row 55
$ws.Range("A55").Value = '4AF08797'
$ws.Range("B55").Value = 'SOL DE JANEIRO CHERIOSA 62 BODY MIST 90 ML'
$ws.Range("C55").Value = 'VARIOS'
$ws.Range("D55").Value = 'Tiene PT'
$ws.Range("E55").Value = 'Tiene ES'
$ws.Range("F55").Value = 'Tiene IT'
$ws.Range("G55").NumberFormat = "@"
$ws.Range("G55").Value = '1'
$ws.Range("H55").Value = 'UND'
$ws.Range("J55").Value = 'Solo Revisión'

# Row 56
$ws.Range("A56").Value = '4AF08798'
$ws.Range("B56").Value = 'SOL DE JANEIRO CHERIOSA 68 BODY MIST 90 ML'
$ws.Range("C56").Value = 'VARIOS'
$ws.Range("D56").Value = 'Tiene PT'
$ws.Range("E56").Value = 'Tiene ES'
$ws.Range("F56").Value = 'Tiene IT'
$ws.Range("G56").NumberFormat = "@"
$ws.Range("G56").Value = '1'
$ws.Range("H56").Value = 'UND'
$ws.Range("J56").Value = 'Solo Revisión'

# Row 57
$ws.Range("A57").Value = '4AF08799'
$ws.Range("B57").Value = 'SOL DE JANEIRO CHERIOSA 68 BODY MIST 240 ML'
$ws.Range("C57").Value = 'VARIOS'
$ws.Range("D57").Value = 'Tiene PT'
$ws.Range("E57").Value = 'Tiene ES'
$ws.Range("F57").Value = 'Tiene IT'
$ws.Range("G57").NumberFormat = "@"
$ws.Range("G57").Value = '1'
$ws.Range("H57").Value = 'UND'
$ws.Range("J57").Value = 'Solo Revisión'

# Row 58
$ws.Range("A58").Value = '4AF08800'
$ws.Range("B58").Value = 'SOL DE JANEIRO RIO RADIANT BODY MIST 90 ML'
$ws.Range("C58").Value = 'VARIOS'
$ws.Range("D58").Value = 'Tiene PT'
$ws.Range("E58").Value = 'Tiene ES'
$ws.Range("F58").Value = 'Tiene IT'
$ws.Range("G58").NumberFormat = "@"
$ws.Range("G58").Value = '1'
$ws.Range("H58").Value = 'UND'
$ws.Range("J58").Value = 'Solo Revisión'

# Row 59
$ws.Range("A59").Value = '6VA39715'
$ws.Range("B59").Value = 'WOOOW KIDNESS BRINGS HAPPINESS ESTUCHE'
$ws.Range("C59").Value = 'VARIOS'
$ws.Range("D59").Value = 'Tiene PT'
$ws.Range("E59").Value = 'Tiene ES'
$ws.Range("F59").Value = 'Tiene IT'
$ws.Range("G59").NumberFormat = "@"
$ws.Range("G59").Value = '1'
$ws.Range("H59").Value = 'UND'
$ws.Range("J59").Value = 'Solo Revisión'

# Row 60
$ws.Range("A60").Value = '1SH00472'
$ws.Range("B60").Value = 'GILLETTE KING C RECORTADORA BARBA STYLE MASTER'
$ws.Range("C60").Value = 'VARIOS'
$ws.Range("D60").Value = 'Tiene PT'
$ws.Range("E60").Value = 'Tiene ES'
$ws.Range("F60").Value = 'Tiene IT'
$ws.Range("G60").NumberFormat = "@"
$ws.Range("G60").Value = '1'
$ws.Range("H60").Value = 'UND'
$ws.Range("J60").Value = 'Solo Revisión'

# Row 61
$ws.Range("A61").Value = '6VA37481'
$ws.Range("B61").Value = 'WOOOOW ESPEJO CUADRADO'
$ws.Range("C61").Value = 'VARIOS'
$ws.Range("D61").Value = 'Tiene PT'
$ws.Range("E61").Value = 'Tiene ES'
$ws.Range("F61").Value = 'Tiene IT'
$ws.Range("G61").NumberFormat = "@"
$ws.Range("G61").Value = '1'
$ws.Range("H61").Value = 'UND'
$ws.Range("J61").Value = 'Solo Revisión'

# Row 62
$ws.Range("A62").Value = '6VA37480'
$ws.Range("B62").Value = 'WOOOOW CORTAUÑAS'
$ws.Range("C62").Value = 'VARIOS'
$ws.Range("D62").Value = 'Tiene PT'
$ws.Range("E62").Value = 'Tiene ES'
$ws.Range("F62").Value = 'Tiene IT'
$ws.Range("G62").NumberFormat = "@"
$ws.Range("G62").Value = '1'
$ws.Range("H62").Value = 'UND'
$ws.Range("J62").Value = 'Revisado y Traducido'

# Row 63
$ws.Range("A63").Value = '6VA37482'
$ws.Range("B63").Value = 'WOOOOW CEPILLO PARA UÑAS'
$ws.Range("C63").Value = 'VARIOS'
$ws.Range("D63").Value = 'Tiene PT'
$ws.Range("E63").Value = 'Tiene ES'
$ws.Range("F63").Value = 'Tiene IT'
$ws.Range("G63").NumberFormat = "@"
$ws.Range("G63").Value = '1'
$ws.Range("H63").Value = 'UND'
$ws.Range("J63").Value = 'Revisado y Traducido'

# Row 64
$ws.Range("A64").Value = '6VA38049'
$ws.Range("B64").Value = 'WOOOOW CEPILLO DOBLE CEJAS & PESTAÑAS'
$ws.Range("C64").Value = 'VARIOS'
$ws.Range("D64").Value = 'Tiene PT'
$ws.Range("E64").Value = 'Tiene ES'
$ws.Range("F64").Value = 'Tiene IT'
$ws.Range("G64").NumberFormat = "@"
$ws.Range("G64").Value = '1'
$ws.Range("H64").Value = 'UND'
$ws.Range("J64").Value = 'Revisado y Traducido'

# Row 65
$ws.Range("A65").Value = '6VA38050'
$ws.Range("B65").Value = 'WOOOOW PEINE METAL PESTAÑAS'
$ws.Range("C65").Value = 'VARIOS'
$ws.Range("D65").Value = 'Tiene PT'
$ws.Range("E65").Value = 'Tiene ES'
$ws.Range("F65").Value = 'Tiene IT'
$ws.Range("G65").NumberFormat = "@"
$ws.Range("G65").Value = '1'
$ws.Range("H65").Value = 'UND'
$ws.Range("J65").Value = 'Revisado y Traducido'

# Row 66
$ws.Range("A66").Value = '6VA37479'
$ws.Range("B66").Value = 'WOOOOW CORTADOR CUTICULA'
$ws.Range("C66").Value = 'VARIOS'
$ws.Range("D66").Value = 'Tiene PT'
$ws.Range("E66").Value = 'Tiene ES'
$ws.Range("F66").Value = 'Tiene IT'
$ws.Range("G66").NumberFormat = "@"
$ws.Range("G66").Value = '1'
$ws.Range("H66").Value = 'UND'
$ws.Range("J66").Value = 'Revisado y Traducido'

# Row 67
$ws.Range("A67").Value = '6VA24971'
$ws.Range("B67").Value = 'ORAL B COMPLETE CEPILLO DENTAL 2UNDS'
$ws.Range("C67").Value = 'VARIOS'
$ws.Range("D67").Value = 'Tiene PT'
$ws.Range("E67").Value = 'Tiene ES'
$ws.Range("F67").Value = 'Tiene IT'
$ws.Range("G67").NumberFormat = "@"
$ws.Range("G67").Value = '2'
$ws.Range("H67").Value = 'UND'
$ws.Range("J67").Value = 'Solo Revisión'

# Row 68
$ws.Range("A68").Value = '6VA37390'
$ws.Range("B68").Value = 'WOOOOW SET BROCHAS PEARL PINK'
$ws.Range("C68").Value = 'VARIOS'
$ws.Range("D68").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E68").Value = 'Tiene ES'
$ws.Range("F68").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G68").NumberFormat = "@"
$ws.Range("G68").Value = '5'
$ws.Range("H68").Value = 'UND'
$ws.Range("J68").Value = 'Revisado y Traducido'

# Row 69
$ws.Range("A69").Value = '0TF27159'
$ws.Range("B69").Value = 'NIVEA Q10 FLUID SPF50 40ML'
$ws.Range("C69").Value = 'TRATAMIENTO CUERPO MANOS'
$ws.Range("D69").Value = 'No Tiene PT - TRADUZIDO'
$ws.Range("E69").Value = 'Tiene ES'
$ws.Range("F69").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G69").NumberFormat = "@"
$ws.Range("G69").Value = '40'
$ws.Range("H69").Value = 'ML'
$ws.Range("J69").Value = 'Revisado y Traducido'

# Row 70
$ws.Range("A70").Value = '0TF27157'
$ws.Range("B70").Value = 'LIPOSAN MANGO 4,8GR'
$ws.Range("C70").Value = 'TRATAMIENTO CUERPO MANOS'
$ws.Range("D70").Value = 'Tiene PT'
$ws.Range("E70").Value = 'Tiene ES'
$ws.Range("F70").Value = 'No Tiene IT - TRADOTTO'
$ws.Range("G70").NumberFormat = "@"
$ws.Range("G70").Value = '4.8'
$ws.Range("H70").Value = 'GR'
$ws.Range("J70").Value = 'Solo Revisión'

# Row 71
$ws.Range("A71").Value = '0TF27144'
$ws.Range("B71").Value = 'SENCE LIP BALM ALOE VERA 2X4,3GRS'
$ws.Range("C71").Value = 'TRATAMIENTO CUERPO MANOS'
$ws.Range("D71").Value = 'Tiene PT'
$ws.Range("E71").Value = 'Tiene ES'
$ws.Range("F71").Value = 'Tiene IT'
$ws.Range("G71").NumberFormat = "@"
$ws.Range("G71").Value = '8.6'
$ws.Range("H71").Value = 'GR'
$ws.Range("J71").Value = 'Solo Revisión'

# Row 72
$ws.Range("A72").Value = '2BG03204'
$ws.Range("B72").Value = 'SENCE GLOW RODILLO FACIAL GIRLS'
$ws.Range("C72").Value = 'BAÑO GEL'
$ws.Range("D72").Value = 'Tiene PT'
$ws.Range("E72").Value = 'Tiene ES'
$ws.Range("F72").Value = 'Tiene IT'
$ws.Range("G72").NumberFormat = "@"
$ws.Range("G72").Value = '1'
$ws.Range("H72").Value = 'UND'
$ws.Range("J72").Value = 'Solo Revisión'
